$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 57
$ws.Range("D2").Value = 211
$ws.Range("C3").Value = 231
$ws.Range("D3").Value = 211
$ws.Range("G3").Value = 1024
$ws.Range("C4").Value = 231
$ws.Range("D4").Value = 278
$ws.Range("G4").Value = 1024
$ws.Range("H4").Value = 768
$ws.Range("C5").Value = 57
$ws.Range("D5").Value = 278
$ws.Range("H5").Value = 768
